$d = $word.ActiveDocument

# Locate the paragraph holding the heuristic-comparison sentence.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*num*incorrect*heuristic the algorithm visited*") {
        $target = $p
        break
    }
}

$pstart = $target.Range.Start
$fullText = $target.Range.Text

# Offsets (relative to the paragraph start) of the single placeholder
# characters "x" and "Y" within "...visited x states, compared to Y with...".
$xOffset = $fullText.IndexOf(" x states") + 1
$yOffset = $fullText.IndexOf(" Y with") + 1

# --- Stage 1: swap in the real numbers, leaving run-splitting for later ---
# (Doing every text assignment before any formatting "touch" keeps the
# engine from re-coalescing earlier splits while we still have more edits
# to make.)
$xStart = $pstart + $xOffset
$xEnd = $xStart + 1
$rX = $d.Range($xStart, $xEnd)
$rX.Text = "1816"

# "Y" is after "x" in the paragraph, and "x" -> "1816" grew the text by 3
# characters, so shift the previously-computed offset accordingly.
$yStart = $pstart + $yOffset + 3
$yEnd = $yStart + 1
$rY = $d.Range($yStart, $yEnd)
$rY.Text = "218"

# --- Stage 2: materialize each newly-typed number as its own run, just
# like Word splits a run at a live edit point even when the formatting is
# unchanged, by touching and reverting a direct character property. ---
$rXNew = $d.Range($xStart, $xStart + 4)
$rXNew.Bold = $true
$rXNew.Bold = $false

$rYNew = $d.Range($yStart, $yStart + 3)
$rYNew.Bold = $true
$rYNew.Bold = $false
